# Automatische test-sync: 2025-06-29 14:48:50
# Adds a new log row (#19) to the "Logs" sheet and updates the
# "Dashboard" category-count summary to reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet: append row 19 with the new test-mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A19").Value2 = "Wil je dit artikel voor me inkopen?"
$logs.Range("B19").Value2 = "mailmind.test@zohomail.eu"
$logs.Range("C19").Value2 = "Testmail #4: Wil je dit artikel voor me inkopen?"
$logs.Range("D19").Value2 = "Bestelling / Levering"
$logs.Range("E19").Value2 = "Beste klant,`nBedankt voor je interesse in ons artikel. Helaas kan ik je op basis van dit bericht niet verder helpen. Kun je meer details geven over welk artikel je wilt inkopen en op welke manier? Zo kan ik je beter assisteren.`nMet vriendelijke groet,`n[Jouw naam]`nE-mailassistent"
$logs.Range("F19").Value2 = "2025-06-29 14:48:22"
$logs.Range("G19").Value2 = "Ja"
$logs.Range("H19").Value2 = "Ja"
$logs.Range("I19").Value2 = "Nee"

# Multi-line content in E19 causes the engine to auto-size the row;
# AutoFit() restores the standard (non-custom) row height so the row
# matches the sheet's normal formatting.
$logs.Rows.Item(19).AutoFit()

# Extend the conditional-formatting ranges so the new row is covered,
# same as Excel normally does when the data region grows. Re-pointing
# each existing rule (instead of deleting/recreating) keeps the rule
# definitions (type/operator/dxf/priority) untouched.
foreach ($col in @("D", "G", "H", "I")) {
    $oldRange = $logs.Range($col + "2:" + $col + "18")
    $newRange = $logs.Range($col + "2:" + $col + "19")
    $rules = $oldRange.FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. "Dashboard" sheet: the new row belongs to category
#    "Bestelling / Levering", so its tally goes from 5 to 6 and the
#    summary rows 2/3 swap places.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value2 = "Bestelling / Levering"
$dash.Range("B2").Value2 = 6
$dash.Range("A3").Value2 = "Openingstijden / Locatie"
$dash.Range("B3").Value2 = 5
